# Update site for cm003 content
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the topic for cm003 (row 4) and mark it as linked
$ws.Range("D4").Value = "Data transformation and exploratory data analysis"
$ws.Range("C4").Value = $true

# Update the active selection shown in the sheet view
$ws.Range("C5").Select()
